$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the test-case data rows (NroSiniestro / NroAnulacion pairs)
$ws.Range("B5").Value = "'0420172010228  "
$ws.Range("C5").Value = "'2027965"

$ws.Range("B6").Value = "'1120170200969 "
$ws.Range("C6").Value = "'0200224 "

$ws.Range("B7").Value = "'1220170301466 "
$ws.Range("C7").Value = "'0300204  "

# Set up page setup (adds pageSetup element on save)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
